$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Transition_Name_Annot").Delete()
$wb.Worksheets.Item("ISTD_Annot").Delete()

$wb.Worksheets.Item("Sample_Annot").Select()
